$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new worker/period data row (row 18) below the existing
#    data rows (16,17), copying formatting from row 17, then fill in
#    the new period's data (new shared string "2509").
$ws.Rows.Item(18).Insert()
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "49606738"
$ws.Range("D18").Value = "CIRINA ELENE CAMARGO ROMERO"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# 2) Update the "VALOR MORA" total and the "Cant. Periodos" count to
#    reflect the newly added period.
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3
